$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2827.4102
$ws.Range("I76").Value = 2463.48
$ws.Range("K76").Value = 2463.48
$ws.Range("M76").Value = -2148.48

$ws.Range("H79").Value = 2827.4102
$ws.Range("I79").Value = 2463.48
$ws.Range("K79").Value = 2463.48
$ws.Range("M79").Value = -1371.48

$ws.Range("H112").Value = 2353
$ws.Range("J112").Value = 3472.75
$ws.Range("L112").Value = 10418.25
$ws.Range("N112").Value = -12634.25

$ws.Range("H132").Value = 4694.759
$ws.Range("I132").Value = 4381.2607
$ws.Range("J132").Value = 5896.5
$ws.Range("K132").Value = 13143.7821
$ws.Range("L132").Value = 17689.5
$ws.Range("M132").Value = -10613.7821
$ws.Range("N132").Value = -22749.5

$ws.Range("H137").Value = 41205.58
$ws.Range("I137").Value = 1856.4
$ws.Range("J137").Value = 94863.55
$ws.Range("K137").Value = 5569.200000000001
$ws.Range("L137").Value = 284590.65
$ws.Range("M137").Value = -3019.200000000001
$ws.Range("N137").Value = -289690.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10756840
$ws.Range("I32").Value = 12347983
$ws.Range("K32").Value = 12347983
$ws.Range("M32").Value = -12347696

$ws.Range("H61").Value = 2143.0417
$ws.Range("I61").Value = 1336.4706
$ws.Range("J61").Value = 4101.857
$ws.Range("K61").Value = 1336.4706
$ws.Range("L61").Value = 4101.857
$ws.Range("M61").Value = -1124.4706
$ws.Range("N61").Value = -4525.857

$ws.Range("H74").Value = 2485.4666
$ws.Range("I74").Value = 2908.913
$ws.Range("J74").Value = 1094.1428
$ws.Range("K74").Value = 2908.913
$ws.Range("L74").Value = 1094.1428
$ws.Range("M74").Value = -2034.913
$ws.Range("N74").Value = -2842.1428

$ws.Range("H77").Value = 2485.4666
$ws.Range("I77").Value = 2908.913
$ws.Range("J77").Value = 1094.1428
$ws.Range("K77").Value = 14544.565
$ws.Range("L77").Value = 5470.714
$ws.Range("M77").Value = -10176.565
$ws.Range("N77").Value = -14206.714

$ws.Range("H136").Value = 2143.0417
$ws.Range("I136").Value = 1336.4706
$ws.Range("J136").Value = 4101.857
$ws.Range("K136").Value = 4009.4118
$ws.Range("L136").Value = 12305.571
$ws.Range("M136").Value = -1459.4118
$ws.Range("N136").Value = -17405.571

$ws.Range("H138").Value = 26315
$ws.Range("J138").Value = 26315
$ws.Range("L138").Value = 26315
$ws.Range("N138").Value = -36595

$ws.Range("H139").Value = 29325
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2381.6667
$ws.Range("I86").Value = 2268
$ws.Range("J86").Value = 2950
$ws.Range("K86").Value = 2268
$ws.Range("L86").Value = 2950
$ws.Range("M86").Value = -1145
$ws.Range("N86").Value = -5196

$ws.Range("H89").Value = 2381.6667
$ws.Range("I89").Value = 2268
$ws.Range("J89").Value = 2950
$ws.Range("K89").Value = 11340
$ws.Range("L89").Value = 14750
$ws.Range("M89").Value = -5724
$ws.Range("N89").Value = -25982

$ws.Range("H99").Value = 2656.3076
$ws.Range("I99").Value = 1462
$ws.Range("J99").Value = 3402.75
$ws.Range("K99").Value = 1462
$ws.Range("L99").Value = 3402.75
$ws.Range("M99").Value = 36
$ws.Range("N99").Value = -6398.75

$ws.Range("H134").Value = 1649.7322
$ws.Range("I134").Value = 1420.8478
$ws.Range("J134").Value = 2702.6
$ws.Range("K134").Value = 4262.5434
$ws.Range("L134").Value = 8107.799999999999
$ws.Range("M134").Value = -1727.5434
$ws.Range("N134").Value = -13177.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2935.5305
$ws.Range("I31").Value = 1696.4286
$ws.Range("J31").Value = 4587.6665
$ws.Range("K31").Value = 1696.4286
$ws.Range("L31").Value = 4587.6665
$ws.Range("M31").Value = -1401.4286
$ws.Range("N31").Value = -5177.6665

$ws.Range("H34").Value = 2935.5305
$ws.Range("I34").Value = 1696.4286
$ws.Range("J34").Value = 4587.6665
$ws.Range("K34").Value = 1696.4286
$ws.Range("L34").Value = 4587.6665
$ws.Range("M34").Value = -1494.4286
$ws.Range("N34").Value = -4991.6665

$ws.Range("H58").Value = 2086.7576
$ws.Range("I58").Value = 1274.826
$ws.Range("J58").Value = 3954.2
$ws.Range("K58").Value = 1274.826
$ws.Range("L58").Value = 3954.2
$ws.Range("M58").Value = -1071.826
$ws.Range("N58").Value = -4360.2

$ws.Range("H132").Value = 2246
$ws.Range("I132").Value = 2272.3635
$ws.Range("J132").Value = 2217
$ws.Range("K132").Value = 6817.0905
$ws.Range("L132").Value = 6651
$ws.Range("M132").Value = -4287.0905
$ws.Range("N132").Value = -11711

$ws.Range("H134").Value = 2466.7646
$ws.Range("I134").Value = 1480.12
$ws.Range("J134").Value = 5207.4443
$ws.Range("K134").Value = 4440.36
$ws.Range("L134").Value = 15622.3329
$ws.Range("M134").Value = -1905.36
$ws.Range("N134").Value = -20692.3329

$ws.Range("H136").Value = 2086.7576
$ws.Range("I136").Value = 1274.826
$ws.Range("J136").Value = 3954.2
$ws.Range("K136").Value = 3824.478
$ws.Range("L136").Value = 11862.6
$ws.Range("M136").Value = -1274.478
$ws.Range("N136").Value = -16962.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 633.1923
$ws.Range("I5").Value = 571.13635
$ws.Range("K5").Value = 1713.40905
$ws.Range("M5").Value = -1601.40905

$ws.Range("H113").Value = 1326696.9
$ws.Range("I113").Value = 2873930.8
$ws.Range("J113").Value = 496.35715
$ws.Range("K113").Value = 8621792.399999999
$ws.Range("L113").Value = 1489.07145
$ws.Range("M113").Value = -8619622.399999999
$ws.Range("N113").Value = -5829.071449999999

$ws.Range("H129").Value = 15437.134
$ws.Range("I129").Value = 2318.7778
$ws.Range("J129").Value = 35114.668
$ws.Range("K129").Value = 6956.3334
$ws.Range("L129").Value = 105344.004
$ws.Range("M129").Value = -1956.3334
$ws.Range("N129").Value = -115344.004

$ws.Range("H130").Value = 1406.7778
$ws.Range("I130").Value = 1207.625
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 3622.875
$ws.Range("L130").Value = 9000
$ws.Range("M130").Value = 1397.125
$ws.Range("N130").Value = -19040

$ws.Range("H135").Value = 633.1923
$ws.Range("I135").Value = 571.13635
$ws.Range("K135").Value = 5140.22715
$ws.Range("M135").Value = -2605.22715

$ws.Range("H136").Value = 4276.2
$ws.Range("I136").Value = 3171.5
$ws.Range("J136").Value = 5933.25
$ws.Range("K136").Value = 9514.5
$ws.Range("L136").Value = 17799.75
$ws.Range("M136").Value = -4414.5
$ws.Range("N136").Value = -27999.75

$ws.Range("H139").Value = 3534.7827
$ws.Range("I139").Value = 885
$ws.Range("J139").Value = 5573.077
$ws.Range("K139").Value = 2655
$ws.Range("L139").Value = 16719.231
$ws.Range("M139").Value = 2485
$ws.Range("N139").Value = -26999.231

$ws.Range("H140").Value = 3718599.2
$ws.Range("I140").Value = 5281588.5
$ws.Range("J140").Value = 6499.875
$ws.Range("K140").Value = 15844765.5
$ws.Range("L140").Value = 19499.625
$ws.Range("M140").Value = -15839585.5
$ws.Range("N140").Value = -29859.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 48.272728
$ws.Range("I2").Value = 48.875
$ws.Range("J2").Value = 46.666668
$ws.Range("K2").Value = 48.875
$ws.Range("L2").Value = 46.666668
$ws.Range("M2").Value = 64.125
$ws.Range("N2").Value = -272.666668

$ws.Range("H141").Value = 70429
$ws.Range("J141").Value = 70429
$ws.Range("L141").Value = 70429
$ws.Range("N141").Value = -80789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1557.71
$ws.Range("I132").Value = 1598.7957
$ws.Range("J132").Value = 1011.8571
$ws.Range("K132").Value = 4796.3871
$ws.Range("L132").Value = 3035.5713
$ws.Range("M132").Value = -2266.3871
$ws.Range("N132").Value = -8095.5713

$ws.Range("H136").Value = 2022.8154
$ws.Range("I136").Value = 1588.8125
$ws.Range("J136").Value = 3248.2354
$ws.Range("K136").Value = 4766.4375
$ws.Range("L136").Value = 9744.706200000001
$ws.Range("M136").Value = -2216.4375
$ws.Range("N136").Value = -14844.7062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45429
$ws.Range("J46").Value = 45429
$ws.Range("L46").Value = 45429
$ws.Range("N46").Value = -45891

$ws.Range("H132").Value = 2380.2856
$ws.Range("I132").Value = 2212.697
$ws.Range("K132").Value = 6638.091
$ws.Range("M132").Value = -4108.091

$ws.Range("H134").Value = 45429
$ws.Range("J134").Value = 45429
$ws.Range("L134").Value = 136287
$ws.Range("N134").Value = -141357

$ws.Range("H136").Value = 2625.2188
$ws.Range("I136").Value = 2573.1372
$ws.Range("J136").Value = 2829.5386
$ws.Range("K136").Value = 7719.4116
$ws.Range("L136").Value = 8488.6158
$ws.Range("M136").Value = -5169.4116
$ws.Range("N136").Value = -13588.6158
